{"js": "// Fix a bug where email body is missing when ML amendment application\n// approved: the Jinja/Django template tag `{% licenced_vessel %}` is\n// missing its `if` keyword in four places in the \"Licensed Vessel\"\n// section (Registration number, Vessel Name, Registered length, Draft).\n// Replace each occurrence of the literal text \"{% licenced_vessel %}\"\n// with \"{% if licenced_vessel %}\", preserving the surrounding run\n// formatting (search ranges carry the formatting of the text they\n// matched, so inserting replacement text with Word.InsertLocation.replace\n// keeps the same font/size/color).\n\nconst searchResults = context.document.body.search(\"{% licenced_vessel %}\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"{% if licenced_vessel %}\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Fix a bug where email body is missing when ML amendment application\n# approved: the Jinja/Django template tag `{% licenced_vessel %}` is\n# missing its `if` keyword in four places in the \"Licensed Vessel\"\n# section (Registration number, Vessel Name, Registered length, Draft).\n# Replace every occurrence of the literal text \"{% licenced_vessel %}\"\n# with \"{% if licenced_vessel %}\" document-wide, preserving the\n# surrounding run formatting (Find/Replace keeps the formatting of the\n# matched text).\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"{% licenced_vessel %}\"\n$find.Replacement.Text = \"{% if licenced_vessel %}\"\n\n$find.Execute(\n    $find.Text,        # FindText\n    $false,             # MatchCase\n    $false,             # MatchWholeWord\n    $false,             # MatchWildcards\n    $false,             # MatchSoundsLike\n    $false,             # MatchAllWordForms\n    $true,              # Forward\n    1,                  # Wrap (wdFindContinue)\n    $false,             # Format\n    $find.Replacement.Text,  # ReplaceWith\n    2                   # Replace (wdReplaceAll)\n) | Out-Null\n"}
